# Applies the 22 Jul 2024 cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.660.59"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.493.42"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'598.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'180.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.14%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.60%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.494.49"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.77%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'7.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.44%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.101.38"
$ws.Range("E13").Value = "  -0.07%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'32.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.45%  "

# Row 15 - TRON
$ws.Range("D15").Value = "'0.134"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "67.634.40"

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.0000178"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.501.36"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.10%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.12%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'392.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'7.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'73.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24 - Polygon
$ws.Range("D24").Value = "'0.543"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.42%  "

# Row 25 - Dai
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  +1.25%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +1.51%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'10.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.34%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -2.05%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  +0.21%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +0.33%  "

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "'2.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'23.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.39%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "'7.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "

# Row 36 - USDe
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +0.55%  "

# Row 38 - Monero
$ws.Range("D38").Value = "'162.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39 - Mantle
$ws.Range("D39").Value = "'0.883"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  +11.95%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "'1.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.55%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "'6.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "

# Row 43 - Filecoin
$ws.Range("D43").Value = "'4.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44 - EnergySwap -> Maker (row swap)
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.832.66"
$ws.Range("E44").Value = "  -0.17%  "

# Row 45 - Maker -> EnergySwap (row swap)
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'26.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.01%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'26.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "

# Row 47 - Hedera
$ws.Range("D47").Value = "'0.0728"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "

# Row 48 - OKB
$ws.Range("D48").Value = "'41.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.21%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +0.20%  "

# Row 50 - Bittensor
$ws.Range("D50").Value = "'336.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.09%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  -0.90%  "
